$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Corrige a quantidade do componente RESISTOR 1000 (linha 28) de 2 para 1.
# As formulas dependentes (Subtotal da linha e Total geral) serao
# recalculadas automaticamente pelo Excel.
$ws.Range("B28").Value = 1
